$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E4").Value = "UNIQUE"
$ws.Range("E7").Value = "UNIQUE"
